# Revert "Drop in files from RMI script"
# Re-introduce the "Texas Data" notes sheet (positioned between "IEA Data"
# and "HPEbP") and fix the natural-gas-reforming efficiency formula on the
# HPEbP sheet, which had incorrectly folded waste heat into the energy
# balance denominator.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Texas Data" worksheet right before "HPEbP" so the
#    final sheet order is: About, IEA Data, Texas Data, HPEbP.
# ------------------------------------------------------------------
$hpebp = $wb.Worksheets.Item("HPEbP")
$texas = $wb.Worksheets.Add($hpebp)
$texas.Name = "Texas Data"

# NOTE: `Worksheets.Add(Before)` inserts the new sheet at the position the
# `$hpebp` handle pointed to, which leaves that variable aliased to the
# newly inserted sheet instead of the original HPEbP sheet. Re-resolve the
# HPEbP worksheet by name so subsequent edits land on the right sheet.
$hpebp = $wb.Worksheets.Item("HPEbP")

# ------------------------------------------------------------------
# 2. Populate the notes left by the author explaining why the Texas
#    numbers mirror the national ones, plus the description of the
#    waste-heat bug found in the NREL-derived calculations.
# ------------------------------------------------------------------
$texas.Range("A1").Value = "There is no reason that these number should be different for Texas."
$texas.Range("A3").Value = "However, I did find an error in their calculations. "
$texas.Range("A5").Value = "They were included waste heat as an energy balance input."
$texas.Range("B6").Value = "for example, page 228 of the NREL report shows gas production as 162 kBtu gas + 2 kBtu electricity = 118 kBtu hydrogen + 46 kBtu waste heat"
$texas.Range("B7").Value = "so, the efficiency (output hydrogen energy vs input energy) would be 118/(162+2)=72%"
$texas.Range("B8").Value = "previously, this spreadhseet (cell 'HPEbP'B3) was calculating the efficiency as 118/(162+2+46)=56%"
$texas.Range("B10").Value = "the IEA number for natural gas reforming efficiency is 76%, so that's a good check that their initial calculation was wrong. "
$texas.Range("A12").Value = "Their other calculations did not include the same mistake."

# Every cell in the used range carries the accent-colored note font.
$texas.Range("A1:I17").Font.ThemeColor = 9

# ------------------------------------------------------------------
# 3. Fix the HPEbP natural-gas-reforming efficiency formula: drop the
#    46 kBtu of waste heat that should never have been in the
#    denominator (118 / (162+2+46) -> 118 / (162+2)). Downstream cells
#    C3:AI3 reference this cell (directly or via shared formulas) and
#    will recalculate automatically.
# ------------------------------------------------------------------
$hpebp.Range("B3").Formula = "=118/(162+2)"

# ------------------------------------------------------------------
# 4. Restore cursor positions / active sheet to match the saved view
#    state, finishing with HPEbP selected as the active tab.
# ------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("B14").Select()

$iea = $wb.Worksheets.Item("IEA Data")
$iea.Range("E18").Select()

$texas.Range("A13").Select()

$hpebp.Range("C12").Select()
$hpebp.Activate()
